# Automatische test-sync: 2025-07-23 22:17:50
#
# Adds the new "Testmail #2" row (M5-bouten order) to the Logs sheet,
# extends its conditional formatting ranges, adds the matching
# aggregate row to the Dashboard sheet, and updates the bar chart's
# series references to include that new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Logs sheet: append row 12
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Kun je 100 stuks M5-bouten bestellen voor klant Jansen?"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Testmail #2: Kun je 100 stuks M5-bouten bestellen voor klant Jansen?"
$logs.Range("D12").Value = "Bestelling / Levering"
$logs.Range("E12").Value = "Beste klant,`nDank je wel voor je e-mail. Om de bestelling voor 100 stuks M5-bouten voor klant Jansen te plaatsen, hebben we wat meer informatie nodig. Zou je ons kunnen voorzien van de specificaties van de M5-bouten (zoals het materiaal, lengte, kopvorm, enz.) zodat we de bestelling correct kunnen verwerken?`nAls je deze informatie kunt verstrekken, zullen we ervoor zorgen dat de bestelling zo spoedig mogelijk wordt geplaatst.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F12").Value = "2025-07-23 22:17:48"
$logs.Range("G12").Value = "Ja"
$logs.Range("H12").Value = "Nee"
$logs.Range("I12").Value = "Ja"
$logs.Range("J12").Value = "Nee"

# Extend the existing conditional-formatting rules so they cover the new row
# (each applies-to range grows from row 11 to row 12, same rules/priorities).
$colRanges = @("D2:D11", "G2:G11", "H2:H11", "I2:I11", "J2:J11")
foreach ($colRange in $colRanges) {
    $newRange = $colRange -replace "11", "12"
    $fcs = $logs.Range($colRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# ---------------------------------------------------------------
# 2) Dashboard sheet: append aggregate row 4
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------
# 3) Update the chart's series to include the new Dashboard row
# ---------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
